$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Bank acount value for "francesco" (row 7) from 810 to 960
$ws.Range("C7").Value = 960
